$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend header row (row 1) with two new columns: P1=14, Q1=15
# These use style index 1 (same as existing header cells like O1), so copy the
# format from the preceding header cell rather than rebuilding it property by
# property (avoids creating redundant style entries).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Update data rows 2-25: columns I, K, M, O swap between 1 and 2,
# and add new columns P and Q (value 2) for each row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2 (new)
}
